$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new, empty paragraph right after the "09/04/18" paragraph.
#    Using Find/Replace with a "^p" code in the replacement text makes the
#    new paragraph correctly inherit the bold+underline paragraph-mark
#    formatting while leaving no stray empty run behind (unlike
#    Range.InsertParagraphAfter(), which always leaves an empty <w:r/>).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("09/04/18", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "09/04/18^p", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Rewrite the "Présentation du bureau..." paragraph:
#      - drop the trailing period after "...fournie par DB"
#      - append the new continuation sentence right after it
#      - append a brand-new paragraph ("Création du squelette...")
#    A single Find/Replace (using "^p" for the paragraph break) achieves
#    this, and naturally relocates the (hidden) _GoBack bookmark - which
#    sits right after the old final "." - into the newly created trailing
#    paragraph.
# ---------------------------------------------------------------------------
$oldTail = "sur la machine fournie par DB."
$newTail = "sur la machine fournie par DB (problème vis-à-vis du proxy qui ne me permettait pas d’installer tout ce que je souhaitais, je devais le faire depuis chez moi).^p" + `
           "Création du squelette du projet, mise en place des premières fonctionnalités (connexion et création de compte).^p"
$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newTail, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert one more empty paragraph after "Création du squelette...", so
#    the (still-hidden) _GoBack bookmark ends up alone in its own trailing
#    paragraph. Doing this as a *separate* Find/Replace call (rather than
#    stacking two "^p^p" in the single replacement above) avoids leaving a
#    stray empty run behind.
# ---------------------------------------------------------------------------
$creationPara = "Création du squelette du projet, mise en place des premières fonctionnalités (connexion et création de compte)."
$d.Content.Find.Execute($creationPara, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $creationPara + "^p", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) The paragraph built in step 2 currently holds its text as a single
#    run. Split it into two runs at the "DB" / " (" boundary by inserting a
#    paragraph break exactly there and then deleting it again: Word keeps
#    the two halves as distinct runs even after the break disappears. This
#    must be the very last edit touching that paragraph, since any further
#    Find/Replace spanning the boundary would coalesce the runs back
#    together.
# ---------------------------------------------------------------------------
$prefix = "Présentation du bureau, des employés et de l’environnement de travail. " + `
          "Installation et configuration des différents composants sur la machine fournie par DB"

$p = $d.Paragraphs.Item(5)
$splitPos = $p.Range.Start + $prefix.Length

$rSplit = $d.Range($splitPos, $splitPos)
$rSplit.InsertParagraphAfter()
$merge = $d.Range($splitPos, $splitPos + 1)
$merge.Delete()

# ---------------------------------------------------------------------------
# 5) Every paragraph in the document now needs "spacing after = 0".
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs.Item($i).SpaceAfter = 0
}
